$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "#enable#" boolean flag column (B) next to userId (A).
$ws.Range("B1").Value = "#enable#"

$ws.Range("B2").Value = $true
$ws.Range("B3").Value = $true

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = $false

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = $false

$ws.Range("A5").Select() | Out-Null
